{"js": "// The author reduced the top page margin (from 1701 twips \u2248 85.05pt to\n// 1134 twips = 56.7pt) on the document's (only) section so the grading\n// table at the end of the TCC report tabulates/fits better before being\n// passed along to DION.\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < sections.items.length; i++) {\n  const pageSetup = sections.items[i].pageSetup;\n  pageSetup.topMargin = 56.7; // points (1134 twips)\n}\n\nawait context.sync();\n", "ps1": "# The author reduced the top page margin (from 1701 twips \u2248 85.05pt to\n# 1134 twips = 56.7pt) so the grading table at the end of the TCC report\n# tabulates/fits better before being passed along to DION.\n$d = $word.ActiveDocument\n\nforeach ($sec in $d.Sections) {\n    $sec.PageSetup.TopMargin = 56.7\n}\n"}
